# Generate Report for Handback
# Update the timestamp strings recorded in the handback status report.
#
# Overview!G2               "Latest HO Xliff Generate Date" for the .md file
# zh-cn!H2                  "Correspond Handoff Datetime"
# zh-cn!K2                  "Correspond Handback DateTime"
# de-de!K2                  "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 19:15:51"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 19:15:47"
$wsZhCn.Range("K2").Value = "2016-08-31 19:16:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-31 19:16:27"
